$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values
$ws.Range("C17").Value = 2.5
$ws.Range("C18").Value = 6

# Clear the values (but keep formatting/style) for C23 and C24
$ws.Range("C23").ClearContents()
$ws.Range("C24").ClearContents()

# Update the active selection to F16
$ws.Range("F16").Select()
